$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (the existing "Notes" column), shifting
# "Notes" and its formatting one column to the right (to L).
$ws.Columns.Item(11).Insert(-4161)  # -4161 = xlShiftToRight

# The new column K should be narrower than the old "Notes" column - just wide
# enough for the "Sales Person" header.
$ws.Columns.Item(11).ColumnWidth = 19.17

# Give the new column its header text.
$ws.Cells.Item(2, 11).Value = "Sales Person"

# Reflect the cursor position left behind by the edit.
$ws.Range("L10").Select() | Out-Null
